$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "29.752.21"
Set-TextValue $ws.Range("E2") "  -1.40%  "
Set-TextValue $ws.Range("D3") "1.888.40"
Set-TextValue $ws.Range("E3") "  -0.96%  "
Set-TextValue $ws.Range("E4") "  -0.10%  "
Set-TextValue $ws.Range("D5") "0.7507"
Set-TextValue $ws.Range("E5") "  +2.89%  "
Set-TextValue $ws.Range("D6") "239.05"
Set-TextValue $ws.Range("E6") "  -1.48%  "
Set-TextValue $ws.Range("E7") "  -0.05%  "
Set-TextValue $ws.Range("B8") "LidoStakedEther"
Set-TextValue $ws.Range("C8") "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
Set-TextValue $ws.Range("D8") "1.887.84"
Set-TextValue $ws.Range("E8") "  -0.34%  "
Set-TextValue $ws.Range("B9") "Cardano"
Set-TextValue $ws.Range("C9") "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextValue $ws.Range("D9") "0.3028"
Set-TextValue $ws.Range("E9") "  -2.94%  "
Set-TextValue $ws.Range("B10") "Solana"
Set-TextValue $ws.Range("C10") "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue $ws.Range("D10") "25.20"
Set-TextValue $ws.Range("E10") "  -4.51%  "
Set-TextValue $ws.Range("B11") "Dogecoin"
Set-TextValue $ws.Range("C11") "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue $ws.Range("D11") "0.06795"
Set-TextValue $ws.Range("E11") "  -1.32%  "
Set-TextValue $ws.Range("B12") "TRON"
Set-TextValue $ws.Range("C12") "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws.Range("D12") "0.07936"
Set-TextValue $ws.Range("E12") "  -0.09%  "
Set-TextValue $ws.Range("D13") "1.894.98"
Set-TextValue $ws.Range("E13") "  -0.29%  "
Set-TextValue $ws.Range("B14") "Polygon"
Set-TextValue $ws.Range("C14") "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws.Range("D14") "0.7410"
Set-TextValue $ws.Range("E14") "  -4.47%  "
Set-TextValue $ws.Range("B15") "Polkadot"
Set-TextValue $ws.Range("C15") "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D15") "5.137"
Set-TextValue $ws.Range("E15") "  -2.25%  "
Set-TextValue $ws.Range("B16") "Litecoin"
Set-TextValue $ws.Range("C16") "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D16") "90.23"
Set-TextValue $ws.Range("E16") "  -1.06%  "
Set-TextValue $ws.Range("B17") "WrappedBTC"
Set-TextValue $ws.Range("C17") "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue $ws.Range("D17") "29.750.67"
Set-TextValue $ws.Range("E17") "  -1.18%  "
Set-TextValue $ws.Range("B18") "Avalanche"
Set-TextValue $ws.Range("C18") "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue $ws.Range("D18") "13.83"
Set-TextValue $ws.Range("E18") "  -2.12%  "
Set-TextValue $ws.Range("B19") "Uniswap"
Set-TextValue $ws.Range("C19") "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws.Range("D19") "5.896"
Set-TextValue $ws.Range("E19") "  +1.23%  "
Set-TextValue $ws.Range("B20") "BitcoinCash"
Set-TextValue $ws.Range("C20") "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws.Range("D20") "241.61"
Set-TextValue $ws.Range("E20") "  +1.02%  "
Set-TextValue $ws.Range("B21") "ShibaInu"
Set-TextValue $ws.Range("C21") "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws.Range("D21") "0.000007634"
Set-TextValue $ws.Range("E21") "  -1.56%  "
Set-TextValue $ws.Range("B22") "Dai"
Set-TextValue $ws.Range("C22") "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D22") "1.001"
Set-TextValue $ws.Range("E22") "  -0.10%  "
Set-TextValue $ws.Range("B23") "BinanceUSD"
Set-TextValue $ws.Range("C23") "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue $ws.Range("D23") "1.001"
Set-TextValue $ws.Range("E23") "  -0.05%  "
Set-TextValue $ws.Range("B24") "Chainlink"
Set-TextValue $ws.Range("C24") "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D24") "6.871"
Set-TextValue $ws.Range("E24") "  -1.15%  "
Set-TextValue $ws.Range("B25") "Monero"
Set-TextValue $ws.Range("C25") "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D25") "165.74"
Set-TextValue $ws.Range("E25") "  +0.71%  "
Set-TextValue $ws.Range("B26") "Cosmos"
Set-TextValue $ws.Range("C26") "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D26") "9.167"
Set-TextValue $ws.Range("E26") "  -2.02%  "
Set-TextValue $ws.Range("B27") "EthereumClassic"
Set-TextValue $ws.Range("C27") "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D27") "18.59"
Set-TextValue $ws.Range("E27") "  -2.52%  "
Set-TextValue $ws.Range("B28") "Stellar"
Set-TextValue $ws.Range("C28") "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D28") "0.1271"
Set-TextValue $ws.Range("E28") "  -0.15%  "
Set-TextValue $ws.Range("B29") "LidoDAOToken"
Set-TextValue $ws.Range("C29") "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws.Range("D29") "2.010"
Set-TextValue $ws.Range("E29") "  -1.91%  "
Set-TextValue $ws.Range("B30") "Toncoin"
Set-TextValue $ws.Range("C30") "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D30") "1.385"
Set-TextValue $ws.Range("E30") "  +2.02%  "
Set-TextValue $ws.Range("B31") "PancakeSwap"
Set-TextValue $ws.Range("C31") "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D31") "1.512"
Set-TextValue $ws.Range("E31") "  -1.84%  "
Set-TextValue $ws.Range("B32") "Filecoin"
Set-TextValue $ws.Range("C32") "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D32") "4.227"
Set-TextValue $ws.Range("E32") "  -1.51%  "
Set-TextValue $ws.Range("B33") "InternetComputer(DFINITY)"
Set-TextValue $ws.Range("C33") "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D33") "3.996"
Set-TextValue $ws.Range("E33") "  -2.03%  "
Set-TextValue $ws.Range("B34") "Hedera"
Set-TextValue $ws.Range("C34") "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D34") "0.05202"
Set-TextValue $ws.Range("E34") "  +1.52%  "
Set-TextValue $ws.Range("B35") "ARBITRUM"
Set-TextValue $ws.Range("C35") "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D35") "1.246"
Set-TextValue $ws.Range("E35") "  -3.14%  "
Set-TextValue $ws.Range("B36") "ImmutableX"
Set-TextValue $ws.Range("C36") "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D36") "0.7232"
Set-TextValue $ws.Range("E36") "  -1.97%  "
Set-TextValue $ws.Range("B37") "HuobiToken"
Set-TextValue $ws.Range("C37") "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D37") "2.708"
Set-TextValue $ws.Range("E37") "  -1.56%  "
Set-TextValue $ws.Range("B38") "VeChain"
Set-TextValue $ws.Range("C38") "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D38") "0.01901"
Set-TextValue $ws.Range("E38") "  -1.46%  "
Set-TextValue $ws.Range("B39") "MXToken"
Set-TextValue $ws.Range("C39") "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D39") "2.760"
Set-TextValue $ws.Range("E39") "  -0.82%  "
Set-TextValue $ws.Range("B40") "FraxShare"
Set-TextValue $ws.Range("C40") "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D40") "6.120"
Set-TextValue $ws.Range("E40") "  -3.60%  "
Set-TextValue $ws.Range("B41") "TheSandbox"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D41") "0.4367"
Set-TextValue $ws.Range("E41") "  -1.37%  "
Set-TextValue $ws.Range("B42") "Aave"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D42") "71.02"
Set-TextValue $ws.Range("E42") "  -4.83%  "
Set-TextValue $ws.Range("B43") "PaxDollar"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue $ws.Range("D43") "1.001"
Set-TextValue $ws.Range("E43") "  +0.03%  "
Set-TextValue $ws.Range("B44") "RenderToken"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D44") "1.873"
Set-TextValue $ws.Range("E44") "  -2.95%  "
Set-TextValue $ws.Range("B45") "TrustWalletToken"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D45") "0.8255"
Set-TextValue $ws.Range("E45") "  -1.06%  "
Set-TextValue $ws.Range("B46") "Aptos"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D46") "7.558"
Set-TextValue $ws.Range("E46") "  +0.04%  "
Set-TextValue $ws.Range("B47") "Quant"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D47") "99.28"
Set-TextValue $ws.Range("E47") "  -1.74%  "
Set-TextValue $ws.Range("B48") "EnergySwap"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D48") "9.672"
Set-TextValue $ws.Range("E48") "  -0.91%  "
Set-TextValue $ws.Range("B49") "RocketPoolETH"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextValue $ws.Range("D49") "2.040.43"
Set-TextValue $ws.Range("E49") "  -0.47%  "
Set-TextValue $ws.Range("B50") "Elrond"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/omwkOTglq+elrond-egld"
Set-TextValue $ws.Range("D50") "35.82"
Set-TextValue $ws.Range("E50") "  -4.57%  "
Set-TextValue $ws.Range("B51") "Cronos"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D51") "0.05943"
Set-TextValue $ws.Range("E51") "  -0.50%  "

Write-Output "Updated 184 cells"
